# Sprint 6 backlog update: rotate board-member reviewers, rename the
# "Admin portal" section to "Frontend changes" with its own new stories,
# and replace the old "Translations" section with new "Backend"
# authentication stories.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3: board-member viewing story - rotate reviewers ---
$ws.Range("D3").Value = "Omar El Desouki"
$ws.Range("E3").Value = "Mostafa Waleed"
$ws.Range("F3").Value = "Beshoy Raef"

# --- Row 4: add/edit board members story - rotate reviewers ---
$ws.Range("D4").Value = "Omar Khaled"
$ws.Range("E4").Value = "Omar El Desouki"
$ws.Range("F4").Value = "Mostafa Waleed"

# --- Row 5: edit board members story - rotate reviewers ---
$ws.Range("D5").Value = "Omar Abdallah"
$ws.Range("E5").Value = "Omar Khaled"
$ws.Range("F5").Value = "Omar El Desouki"

# --- Row 7: section heading renamed "Admin portal" -> "Frontend changes" ---
$ws.Range("B7").Value = "Frontend changes"

# --- Row 8: admin portal redirect story - rotate reviewers (assignee stays) ---
$ws.Range("D8").Value = "Ahmed Osama"
$ws.Range("E8").Value = "Omar Abdallah"
$ws.Range("F8").Value = "Omar Khaled"

# --- Row 9 (new): edit/delete own comments story ---
# Copy formatting from the row above so the new row matches the rest of
# the table's style, then fill in the values.
$ws.Range("A8:F8").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A9").Value = 2.2
$ws.Range("B9").Value = "As a user, I can edit or delete my comments"
$ws.Range("C9").Value = 3
$ws.Range("D9").Value = "Yosri Khaled"
$ws.Range("E9").Value = "Ahmed Osama"
$ws.Range("F9").Value = "Omar Abdallah"

# --- Row 10: was the "Translations" section header, now a data row ---
# C10:F10 had no prior formatting (row 10 used to be a two-cell section
# header), so pull the style from row 9 before filling in the values.
$ws.Range("C9:F9").Copy()
$ws.Range("C10:F10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A10").Value = 2.3
$ws.Range("B10").Value = "As a user, I can delete my company application"
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = "Abed Hossam"
$ws.Range("E10").Value = "Yosri Khaled"
$ws.Range("F10").Value = "Ahmed Osama"

# --- Row 11: old "Implement translation API" row is removed entirely ---
$ws.Range("A11:F11").Clear()

# --- Row 12: becomes the new "Backend" section header ---
$ws.Range("A12").Value = 3
$ws.Range("B12").Value = "Backend"
$ws.Range("C12:F12").Clear()

# --- Row 13: becomes "Authentication for user controller" ---
$ws.Range("A13").Value = 3.1
$ws.Range("B13").Value = "Authentication for user controller"
$ws.Range("C13").Value = 5
$ws.Range("D13").Value = "Ahmed Ashraf"
$ws.Range("E13").Value = "Abed Hossam"
$ws.Range("F13").Value = "Yosri Khaled"

# --- Row 14 (new): "Authentication for company controller" ---
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A14").Value = 3.2
$ws.Range("B14").Value = "Authentication for company controller"
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = "Beshoy Raef"
$ws.Range("E14").Value = "Ahmed Ashraf"
$ws.Range("F14").Value = "Abed Hossam"

# --- Row 15: becomes "Authentication for other controllers" ---
# C15:F15 had no prior formatting (row 15 used to be the "Backend"
# section header), so pull the style from row 14 before filling in values.
$ws.Range("C14:F14").Copy()
$ws.Range("C15:F15").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A15").Value = 3.3
$ws.Range("B15").Value = "Authentication for other controllers"
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = "Mostafa Waleed"
$ws.Range("E15").Value = "Beshoy Raef"
$ws.Range("F15").Value = "Ahmed Ashraf"

# --- Rows 16-17: old content shifted up / consolidated, now empty ---
$ws.Range("A16:F17").Clear()
